$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3278.5715
$ws.Range("J51").Value = 3450
$ws.Range("L51").Value = 3450
$ws.Range("N51").Value = -4418
$ws.Range("H80").Value = 1140
$ws.Range("I80").Value = 855.25
$ws.Range("J80").Value = 1266.5555
$ws.Range("K80").Value = 2565.75
$ws.Range("L80").Value = 3799.6665
$ws.Range("M80").Value = -1567.75
$ws.Range("N80").Value = -5795.666499999999
$ws.Range("H83").Value = 1140
$ws.Range("I83").Value = 855.25
$ws.Range("J83").Value = 1266.5555
$ws.Range("K83").Value = 7697.25
$ws.Range("L83").Value = 11398.9995
$ws.Range("M83").Value = -2705.25
$ws.Range("N83").Value = -21382.9995
$ws.Range("H87").Value = 58925
$ws.Range("J87").Value = 65628.57000000001
$ws.Range("L87").Value = 65628.57000000001
$ws.Range("N87").Value = -68124.57000000001
$ws.Range("H90").Value = 58925
$ws.Range("J90").Value = 65628.57000000001
$ws.Range("L90").Value = 196885.71
$ws.Range("N90").Value = -209365.71
$ws.Range("H92").Value = 4900
$ws.Range("I92").Value = 4250
$ws.Range("J92").Value = 5333.3335
$ws.Range("K92").Value = 4250
$ws.Range("L92").Value = 5333.3335
$ws.Range("M92").Value = -3002
$ws.Range("N92").Value = -7829.3335
$ws.Range("H138").Value = 2649.3333
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2649.3333
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7947.999899999999
$ws.Range("N138").Value = -18227.9999
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4425.108
$ws.Range("I32").Value = 4425.108
$ws.Range("K32").Value = 4425.108
$ws.Range("M32").Value = -4138.108
$ws.Range("H45").Value = 3127
$ws.Range("I45").Value = 2620
$ws.Range("K45").Value = 2620
$ws.Range("M45").Value = -2243
$ws.Range("H63").Value = 6354.1665
$ws.Range("I63").Value = 1849.5
$ws.Range("J63").Value = 8606.5
$ws.Range("K63").Value = 1849.5
$ws.Range("L63").Value = 8606.5
$ws.Range("M63").Value = -1163.5
$ws.Range("N63").Value = -9978.5
$ws.Range("H66").Value = 6354.1665
$ws.Range("I66").Value = 1849.5
$ws.Range("J66").Value = 8606.5
$ws.Range("K66").Value = 9247.5
$ws.Range("L66").Value = 43032.5
$ws.Range("M66").Value = -5815.5
$ws.Range("N66").Value = -49896.5
$ws.Range("H101").Value = 68351.336
$ws.Range("J101").Value = 68351.336
$ws.Range("L101").Value = 68351.336
$ws.Range("N101").Value = -74841.336
$ws.Range("H110").Value = 4010
$ws.Range("I110").Value = 1400
$ws.Range("K110").Value = 1400
$ws.Range("M110").Value = 645

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H86").Value = 6818.4116
$ws.Range("I86").Value = 5696.3
$ws.Range("J86").Value = 8421.429
$ws.Range("K86").Value = 5696.3
$ws.Range("L86").Value = 8421.429
$ws.Range("M86").Value = -4573.3
$ws.Range("N86").Value = -10667.429
$ws.Range("H89").Value = 6818.4116
$ws.Range("I89").Value = 5696.3
$ws.Range("J89").Value = 8421.429
$ws.Range("K89").Value = 28481.5
$ws.Range("L89").Value = 42107.145
$ws.Range("M89").Value = -22865.5
$ws.Range("N89").Value = -53339.145
$ws.Range("H105").Value = 2285
$ws.Range("I105").Value = 2296.8
$ws.Range("K105").Value = 2296.8
$ws.Range("M105").Value = -549.8000000000002
$ws.Range("H134").Value = 4392.885
$ws.Range("I134").Value = 4168.6
$ws.Range("K134").Value = 12505.8
$ws.Range("M134").Value = -9970.800000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1793.8889
$ws.Range("I16").Value = 1793.8889
$ws.Range("K16").Value = 1793.8889
$ws.Range("M16").Value = -1506.8889
$ws.Range("H31").Value = 5482.575
$ws.Range("J31").Value = 9803.294
$ws.Range("L31").Value = 9803.294
$ws.Range("N31").Value = -10393.294
$ws.Range("H34").Value = 5482.575
$ws.Range("J34").Value = 9803.294
$ws.Range("L34").Value = 9803.294
$ws.Range("N34").Value = -10207.294
$ws.Range("H60").Value = 1000
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -489
$ws.Range("N60").ClearContents()
$ws.Range("H113").Value = 1793.8889
$ws.Range("I113").Value = 1793.8889
$ws.Range("K113").Value = 1793.8889
$ws.Range("M113").Value = 376.1111000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 1899999.4
$ws.Range("I128").Value = 1899999.4
$ws.Range("K128").Value = 5699998.199999999
$ws.Range("M128").Value = -5695018.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H102").Value = 1961
$ws.Range("I102").Value = 1596.36
$ws.Range("K102").Value = 1596.36
$ws.Range("M102").Value = 25.6400000000001
$ws.Range("H113").Value = 9777.111000000001
$ws.Range("I113").Value = 8999.5
$ws.Range("J113").Value = 9999.286
$ws.Range("K113").Value = 8999.5
$ws.Range("L113").Value = 9999.286
$ws.Range("M113").Value = -6829.5
$ws.Range("N113").Value = -14339.286
$ws.Range("H136").Value = 76666.664
$ws.Range("J136").Value = 76666.664
$ws.Range("L136").Value = 229999.992
$ws.Range("N136").Value = -235099.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3351.4583
$ws.Range("I61").Value = 2421.75
$ws.Range("K61").Value = 2421.75
$ws.Range("M61").Value = -2219.75
$ws.Range("H68").Value = 7480.7
$ws.Range("I68").Value = 1936.3334
$ws.Range("J68").Value = 9856.857
$ws.Range("K68").Value = 1936.3334
$ws.Range("L68").Value = 9856.857
$ws.Range("M68").Value = -1187.3334
$ws.Range("N68").Value = -11354.857
$ws.Range("H71").Value = 7480.7
$ws.Range("I71").Value = 1936.3334
$ws.Range("J71").Value = 9856.857
$ws.Range("K71").Value = 9681.666999999999
$ws.Range("L71").Value = 49284.285
$ws.Range("M71").Value = -5937.666999999999
$ws.Range("N71").Value = -56772.285
$ws.Range("H82").Value = 5908.909
$ws.Range("J82").Value = 6983.1665
$ws.Range("L82").Value = 6983.1665
$ws.Range("N82").Value = -7705.1665
$ws.Range("H85").Value = 5908.909
$ws.Range("J85").Value = 6983.1665
$ws.Range("L85").Value = 6983.1665
$ws.Range("N85").Value = -9479.166499999999
$ws.Range("H93").Value = 2435.75
$ws.Range("J93").Value = 2555
$ws.Range("L93").Value = 2555
$ws.Range("N93").Value = -5051
$ws.Range("H100").Value = 5659.9
$ws.Range("I100").Value = 2212.25
$ws.Range("K100").Value = 2212.25
$ws.Range("M100").Value = -1671.25
$ws.Range("H113").Value = 3351.4583
$ws.Range("I113").Value = 2421.75
$ws.Range("K113").Value = 2421.75
$ws.Range("M113").Value = -251.75
$ws.Range("H132").Value = 3797
$ws.Range("I132").Value = 3797
$ws.Range("K132").Value = 11391
$ws.Range("M132").Value = -8861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49297.2
$ws.Range("J46").Value = 49297.2
$ws.Range("L46").Value = 49297.2
$ws.Range("N46").Value = -49759.2
$ws.Range("H134").Value = 49297.2
$ws.Range("J134").Value = 49297.2
$ws.Range("L134").Value = 147891.6
$ws.Range("N134").Value = -152961.6
